$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the stored calendar dates (fechaInicial / FechaFinal) in B2 and
# B3 so the workbook works with any range of calendar dates.
# B2: 2020-01-01 (serial 43831) -> 2020-09-29 (serial 44103)
# B3: 2020-01-16 (serial 43846) -> 2021-01-14 (serial 44210)
# Cells already carry a date number format (style id 1), so writing the
# bare serial number stores a clean whole-day value with no time-of-day
# fraction, exactly like typing the date into the cell in Excel.
$ws.Range("B2").Value = 44103
$ws.Range("B3").Value = 44210

# Reset/maximize the workbook window view (best-effort UI state).
$excel.ActiveWindow.WindowState = -4137
